$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format for the editable data range so that
# numeric-looking price strings (e.g. "251.33") are stored as text, matching
# the original inline-string data rather than being parsed into numbers.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '41.872.52'
$ws.Range("E2").Value = '  +1.15%  '

# Row 3
$ws.Range("D3").Value = '2.205.46'
$ws.Range("E3").Value = '  +0.60%  '

# Row 4
$ws.Range("E4").Value = '  -0.04%  '

# Row 5
$ws.Range("D5").Value = '251.33'
$ws.Range("E5").Value = '  -1.52%  '

# Row 6
$ws.Range("E6").Value = '  -0.59%  '

# Row 7
$ws.Range("D7").Value = '67.46'
$ws.Range("E7").Value = '  -1.01%  '

# Row 8
$ws.Range("E8").Value = '  -0.08%  '

# Row 9
$ws.Range("D9").Value = '0.617'
$ws.Range("E9").Value = '  +6.21%  '

# Row 10
$ws.Range("D10").Value = '38.75'
$ws.Range("E10").Value = '  +2.07%  '

# Row 11
$ws.Range("D11").Value = '59.30'
$ws.Range("E11").Value = '  +1.87%  '

# Row 12
$ws.Range("D12").Value = '0.0936'
$ws.Range("E12").Value = '  -0.82%  '

# Row 13
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").Value = '0.104'
$ws.Range("E13").Value = '  -0.08%  '

# Row 14
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").Value = '6.97'
$ws.Range("E14").Value = '  -1.51%  '

# Row 15
$ws.Range("D15").Value = '2.538.30'
$ws.Range("E15").Value = '  +0.88%  '

# Row 16
$ws.Range("D16").Value = '0.872'
$ws.Range("E16").Value = '  -0.08%  '

# Row 17
$ws.Range("D17").Value = '14.49'
$ws.Range("E17").Value = '  -0.18%  '

# Row 18
$ws.Range("D18").Value = '2.204.86'
$ws.Range("E18").Value = '  +0.95%  '

# Row 19
$ws.Range("D19").Value = '41.811.19'
$ws.Range("E19").Value = '  +1.36%  '

# Row 20
$ws.Range("E20").Value = '  +0.50%  '

# Row 21
$ws.Range("D21").Value = '72.31'
$ws.Range("E21").Value = '  +0.17%  '

# Row 22
$ws.Range("D22").Value = '6.13'
$ws.Range("E22").Value = '  -1.89%  '

# Row 23
$ws.Range("D23").Value = '231.04'
$ws.Range("E23").Value = '  -0.75%  '

# Row 24
$ws.Range("D24").Value = '2.01'
$ws.Range("E24").Value = '  -3.00%  '

# Row 25
$ws.Range("D25").Value = '3.89'
$ws.Range("E25").Value = '  +2.19%  '

# Row 26
$ws.Range("E26").Value = '  +0.19%  '

# Row 27
$ws.Range("D27").Value = '11.13'
$ws.Range("E27").Value = '  -6.67%  '

# Row 28
$ws.Range("D28").Value = '2.41'
$ws.Range("E28").Value = '  -4.72%  '

# Row 29
$ws.Range("D29").Value = '3.68'
$ws.Range("E29").Value = '  -1.41%  '

# Row 30
$ws.Range("E30").Value = '  -1.20%  '

# Row 31
$ws.Range("D31").Value = '166.42'
$ws.Range("E31").Value = '  -1.75%  '

# Row 32
$ws.Range("E32").Value = '  -1.29%  '

# Row 33
$ws.Range("E33").Value = '  +0.98%  '

# Row 34
$ws.Range("D34").Value = '5.86'
$ws.Range("E34").Value = '  +7.16%  '

# Row 35
$ws.Range("D35").Value = '0.0778'
$ws.Range("E35").Value = '  +6.85%  '

# Row 36
$ws.Range("E36").Value = '  -0.33%  '

# Row 37
$ws.Range("D37").Value = '25.91'
$ws.Range("E37").Value = '  +2.49%  '

# Row 38
$ws.Range("D38").Value = '4.58'
$ws.Range("E38").Value = '  -0.44%  '

# Row 39
$ws.Range("E39").Value = '  +1.13%  '

# Row 40
$ws.Range("D40").Value = '0.0308'
$ws.Range("E40").Value = '  +2.95%  '

# Row 41
$ws.Range("E41").Value = '  -0.52%  '

# Row 42
$ws.Range("D42").Value = '5.18'
$ws.Range("E42").Value = '  +6.60%  '

# Row 43
$ws.Range("B43").Value = 'Celestia'
$ws.Range("C43").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D43").Value = '11.99'
$ws.Range("E43").Value = '  -1.88%  '

# Row 44
$ws.Range("B44").Value = 'THORChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D44").Value = '5.62'
$ws.Range("E44").Value = '  -2.38%  '

# Row 45
$ws.Range("B45").Value = 'Algorand'
$ws.Range("C45").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D45").Value = '0.196'
$ws.Range("E45").Value = '  -3.51%  '

# Row 46
$ws.Range("B46").Value = 'MultiversX'
$ws.Range("C46").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D46").Value = '59.69'
$ws.Range("E46").Value = '  -7.39%  '

# Row 47
$ws.Range("E47").Value = '  -0.79%  '

# Row 48
$ws.Range("D48").Value = '0.0995'
$ws.Range("E48").Value = '  -2.01%  '

# Row 49
$ws.Range("E49").Value = '  -0.44%  '

# Row 50
$ws.Range("E50").Value = '  +0.17%  '

# Row 51
$ws.Range("D51").Value = '2.83'
$ws.Range("E51").Value = '  +4.68%  '

# Restore the default (Normal) style on column D so no stray number format
# remains applied to the cells now that the text values are stored.
$priceRange.Style = "Normal"
